$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=2; Col=2; Value=0.2053231939163498}
    @{Row=2; Col=3; Value=0.5627376425855514}
    @{Row=2; Col=10; Value=0.01140684410646388}
    @{Row=2; Col=16; Value=0.1330798479087452}
    @{Row=2; Col=19; Value=0.08745247148288973}
    @{Row=3; Col=2; Value=0.006535947712418301}
    @{Row=3; Col=3; Value=0.0392156862745098}
    @{Row=3; Col=10; Value=0.0261437908496732}
    @{Row=3; Col=16; Value=0.7516339869281046}
    @{Row=3; Col=19; Value=0.1764705882352941}
    @{Row=4; Col=10; Value=0.05128205128205128}
    @{Row=4; Col=16; Value=0.6923076923076923}
    @{Row=4; Col=19; Value=0.2564102564102564}
    @{Row=6; Col=2; Value=0.06967213114754098}
    @{Row=6; Col=4; Value=0.01229508196721311}
    @{Row=6; Col=6; Value=0.09426229508196721}
    @{Row=6; Col=10; Value=0.2663934426229508}
    @{Row=6; Col=15; Value=0.02459016393442623}
    @{Row=6; Col=17; Value=0.1516393442622951}
    @{Row=6; Col=18; Value=0.03688524590163934}
    @{Row=6; Col=19; Value=0.3442622950819672}
    @{Row=7; Col=2; Value=0.08695652173913043}
    @{Row=7; Col=4; Value=0.02173913043478261}
    @{Row=7; Col=5; Value=0.004347826086956522}
    @{Row=7; Col=6; Value=0.07391304347826087}
    @{Row=7; Col=10; Value=0.1652173913043478}
    @{Row=7; Col=15; Value=0.008695652173913044}
    @{Row=7; Col=17; Value=0.208695652173913}
    @{Row=7; Col=18; Value=0.06086956521739131}
    @{Row=7; Col=19; Value=0.3695652173913043}
    @{Row=8; Col=2; Value=0.1006711409395973}
    @{Row=8; Col=4; Value=0.01565995525727069}
    @{Row=8; Col=5; Value=0.002237136465324385}
    @{Row=8; Col=6; Value=0.06935123042505593}
    @{Row=8; Col=10; Value=0.08501118568232663}
    @{Row=8; Col=15; Value=0.02237136465324385}
    @{Row=8; Col=17; Value=0.1968680089485459}
    @{Row=8; Col=18; Value=0.08724832214765101}
    @{Row=8; Col=19; Value=0.4205816554809844}
    @{Row=9; Col=2; Value=0.09865470852017937}
    @{Row=9; Col=4; Value=0.008968609865470852}
    @{Row=9; Col=6; Value=0.08520179372197309}
    @{Row=9; Col=10; Value=0.1121076233183857}
    @{Row=9; Col=15; Value=0.02690582959641256}
    @{Row=9; Col=17; Value=0.2197309417040359}
    @{Row=9; Col=18; Value=0.1121076233183857}
    @{Row=9; Col=19; Value=0.336322869955157}
    @{Row=10; Col=2; Value=0.09751609935602576}
    @{Row=10; Col=4; Value=0.0202391904323827}
    @{Row=10; Col=5; Value=0.002759889604415824}
    @{Row=10; Col=6; Value=0.07543698252069918}
    @{Row=10; Col=10; Value=0.09567617295308188}
    @{Row=10; Col=15; Value=0.01931922723091076}
    @{Row=10; Col=17; Value=0.2207911683532659}
    @{Row=10; Col=18; Value=0.07359705611775529}
    @{Row=10; Col=19; Value=0.3946642134314627}
    @{Row=11; Col=7; Value=0.1505681818181818}
    @{Row=11; Col=10; Value=0.08806818181818182}
    @{Row=11; Col=11; Value=0.2073863636363636}
    @{Row=11; Col=12; Value=0.5397727272727273}
    @{Row=11; Col=19; Value=0.01420454545454545}
    @{Row=12; Col=7; Value=0.745}
    @{Row=12; Col=10; Value=0.175}
    @{Row=12; Col=11; Value=0.005}
    @{Row=12; Col=12; Value=0.05}
    @{Row=12; Col=19; Value=0.025}
    @{Row=13; Col=7; Value=0.7727272727272727}
    @{Row=13; Col=10; Value=0.2045454545454546}
    @{Row=13; Col=19; Value=0.02272727272727273}
    @{Row=15; Col=6; Value=0.03238866396761134}
    @{Row=15; Col=8; Value=0.1659919028340081}
    @{Row=15; Col=9; Value=0.0931174089068826}
    @{Row=15; Col=10; Value=0.3036437246963563}
    @{Row=15; Col=11; Value=0.0931174089068826}
    @{Row=15; Col=13; Value=0.008097165991902834}
    @{Row=15; Col=15; Value=0.04453441295546558}
    @{Row=15; Col=19; Value=0.2591093117408907}
    @{Row=16; Col=6; Value=0.01169590643274854}
    @{Row=16; Col=8; Value=0.2046783625730994}
    @{Row=16; Col=9; Value=0.06432748538011696}
    @{Row=16; Col=10; Value=0.4327485380116959}
    @{Row=16; Col=11; Value=0.08771929824561403}
    @{Row=16; Col=13; Value=0.01754385964912281}
    @{Row=16; Col=15; Value=0.07602339181286549}
    @{Row=16; Col=19; Value=0.1052631578947368}
    @{Row=17; Col=6; Value=0.03017241379310345}
    @{Row=17; Col=8; Value=0.2133620689655172}
    @{Row=17; Col=9; Value=0.08836206896551724}
    @{Row=17; Col=10; Value=0.3275862068965517}
    @{Row=17; Col=11; Value=0.1185344827586207}
    @{Row=17; Col=13; Value=0.02801724137931035}
    @{Row=17; Col=14; Value=0.004310344827586207}
    @{Row=17; Col=15; Value=0.07758620689655173}
    @{Row=17; Col=19; Value=0.1120689655172414}
    @{Row=18; Col=6; Value=0.005988023952095809}
    @{Row=18; Col=8; Value=0.1736526946107785}
    @{Row=18; Col=9; Value=0.125748502994012}
    @{Row=18; Col=10; Value=0.3173652694610778}
    @{Row=18; Col=11; Value=0.1317365269461078}
    @{Row=18; Col=13; Value=0.005988023952095809}
    @{Row=18; Col=15; Value=0.08383233532934131}
    @{Row=18; Col=19; Value=0.155688622754491}
    @{Row=19; Col=6; Value=0.02117263843648208}
    @{Row=19; Col=8; Value=0.1978827361563518}
    @{Row=19; Col=9; Value=0.1034201954397394}
    @{Row=19; Col=10; Value=0.3184039087947882}
    @{Row=19; Col=11; Value=0.1262214983713355}
    @{Row=19; Col=13; Value=0.02361563517915309}
    @{Row=19; Col=14; Value=0.003257328990228013}
    @{Row=19; Col=15; Value=0.07491856677524431}
    @{Row=19; Col=19; Value=0.1311074918566775}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value2 = $u.Value
}

Write-Output "Updated $($updates.Count) cells"